# Apply updated capacity values in column C (planned retire groupings)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    10 = 1.528
    13 = 6.064
    18 = 1430.5294
    21 = 19.691
    34 = 1.77
    37 = 0.32
    42 = 0.426
    45 = 0.824
    58 = 0.3
    61 = 2.9708
    62 = 0.334
    65 = 1.832
    66 = 1146.536
    69 = 3.26
    90 = 0.7869999999999999
    93 = 0.411
    102 = 36.5322
    105 = 46.5096
    126 = 0.122
    129 = 1.907
    130 = 1686.7424
    133 = 163.8495
    134 = 144.307
    137 = 105.3884
    142 = 9.83
    145 = 16.4365
    150 = 1823.5223
    153 = 187.5155
    178 = 0.25
    181 = 0.884
    186 = 50.0116
    189 = 0.195
    202 = 0.336
    205 = 5.83
    210 = 51.7432
    213 = 1.4168
    214 = 12.433
    217 = 1.173
    314 = 201.6404
    317 = 145.3945
    342 = 11.1691
    345 = 7.4809
    354 = 33.7861
    357 = 4.025
    362 = 3.737
    365 = 0.73
    370 = 0.124
    373 = 0
    378 = 9.601000000000001
    381 = 8.4
    382 = 2.0343
    385 = 0.8662000000000001
    430 = 2.172
    433 = 0
    434 = 79.8038
    437 = 76.82729999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $updates[$row]
}

$wb.Save()
